$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row = 13; dev = "Developer12"; tester = "tester12"; name = "SELF PACED ONLINE TRAINING" },
    @{ row = 14; dev = "Developer13"; tester = "tester13"; name = "IN DEPTH MATERIAL" },
    @{ row = 15; dev = "Developer14"; tester = "tester14"; name = "LIFETIME INSTRUCTOR SUPPORT" },
    @{ row = 16; dev = "Developer15"; tester = "tester15"; name = "RESUME PREPARATION" }
)

foreach ($r in $rows) {
    $rowIndex = $r.row
    $name = $r.name

    $ws.Cells.Item($rowIndex, 1).Value = "validate $name"
    $ws.Cells.Item($rowIndex, 2).Value = $r.dev
    $ws.Cells.Item($rowIndex, 3).Value = "landing Page"
    $ws.Cells.Item($rowIndex, 4).Value = "Validate $name"
    $ws.Cells.Item($rowIndex, 5).Value = "1/4"
    $ws.Cells.Item($rowIndex, 6).Value = "Pre-conditions"
    $ws.Cells.Item($rowIndex, 7).Value = "N/A"
    $ws.Cells.Item($rowIndex, 8).Value = "Open Site`nClick in New window`nValidar text $name"
    $ws.Cells.Item($rowIndex, 9).Value = "Text should be in site"
    $ws.Cells.Item($rowIndex, 10).Value = $r.tester
    $ws.Cells.Item($rowIndex, 11).Value = "Text isnt in front of the site"
    $ws.Cells.Item($rowIndex, 12).Value = "Fail"
    $ws.Cells.Item($rowIndex, 13).Value = "Test automation failed"
    $ws.Cells.Item($rowIndex, 14).Value = "Or the output"
}
